# The data table on Sheet1 gets one new weekly price record inserted as
# row 85 (pushing the existing rows 85-161 down to 86-162, preserving
# their data unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 85; everything below (old rows
# 85-161) shifts down to 86-162 with all of its data intact.
$ws.Rows("85:85").Insert()

# Populate the newly inserted row 85 with the new record.
$ws.Cells.Item(85, 1).Value2  = 7
$ws.Cells.Item(85, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(85, 3).Value2  = "Ñuble"
$ws.Cells.Item(85, 4).Value2  = 44554
$ws.Cells.Item(85, 5).Value2  = 16
$ws.Cells.Item(85, 6).Value2  = 100112017
$ws.Cells.Item(85, 7).Value2  = "Apio"
$ws.Cells.Item(85, 8).Value2  = "Americana (o)"
$ws.Cells.Item(85, 9).Value2  = "Primera"
$ws.Cells.Item(85, 10).Value2 = 80
$ws.Cells.Item(85, 11).Value2 = 8000
$ws.Cells.Item(85, 12).Value2 = 8500
$ws.Cells.Item(85, 13).Value2 = 8250
$ws.Cells.Item(85, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(85, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(85, 16).Value2 = 1375
$ws.Cells.Item(85, 17).Value2 = 6
$ws.Cells.Item(85, 18).Value2 = "Hortaliza"

# Make sure the new date cell keeps the same date formatting as the
# rest of column D (style carries over from the Insert, but set it
# explicitly to be safe).
$ws.Cells.Item(85, 4).NumberFormat = $ws.Cells.Item(86, 4).NumberFormat
